# "Fruta / hortaliza, semanal" -- add a new weekly price record.
#
# A new data row is inserted right above the current row 15 (pushing the
# existing rows 15..105 down to 16..106, and growing the sheet's used range
# from A1:T105 to A1:T106). The freshly inserted row is then populated with
# the new Granada ("Sin especificar" / "Primera") price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 15; everything below shifts down one row.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data.
$ws.Cells.Item(15, 1).Value  = 10                                   # Mercado ID
$ws.Cells.Item(15, 2).Value  = "Vega Modelo de Temuco"               # Mercado
$ws.Cells.Item(15, 3).Value  = "La Araucanía"                        # Región
$ws.Cells.Item(15, 4).Value  = 44635                                 # Fecha (2022-03-15)
$ws.Cells.Item(15, 5).Value  = 9                                     # Codreg
$ws.Cells.Item(15, 6).Value  = "Fruta"                                # Tipo
$ws.Cells.Item(15, 7).Value  = 100104                                 # Producto ID
$ws.Cells.Item(15, 8).Value  = "Frutos de pepita"                     # Producto
$ws.Cells.Item(15, 9).Value  = 100104001                              # Categoría ID
$ws.Cells.Item(15, 10).Value = "Granada"                              # Categoría
$ws.Cells.Item(15, 11).Value = "Sin especificar"                      # Variedad
$ws.Cells.Item(15, 12).Value = "Primera"                              # Calidad
$ws.Cells.Item(15, 13).Value = 25                                     # Volumen
$ws.Cells.Item(15, 14).Value = 20000                                  # Precio mínimo
$ws.Cells.Item(15, 15).Value = 20000                                  # Precio máximo
$ws.Cells.Item(15, 16).Value = 20000                                  # Precio promedio ponderado
$ws.Cells.Item(15, 17).Value = "$/bandeja 10 kilos granel"            # Unidad de comercialización
$ws.Cells.Item(15, 18).Value = "Provincia de Limarí"                  # Origen
$ws.Cells.Item(15, 19).Value = 2000                                   # Precio $/Kg
$ws.Cells.Item(15, 20).Value = 10                                     # Kg / unidad
